$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (single-dot decimal values).
# Cells whose new values contain two dots (e.g. "24.391.48") are left alone
# since Excel cannot parse them as numbers and keeps them as text naturally.
$ws.Range("D4:D10").NumberFormat = "@"
$ws.Range("D12:D16").NumberFormat = "@"
$ws.Range("D18:D23").NumberFormat = "@"
$ws.Range("D25:D32").NumberFormat = "@"
$ws.Range("D34:D45").NumberFormat = "@"
$ws.Range("D47:D51").NumberFormat = "@"

# Update Price (D) column
$ws.Range("D2").Value = "24.391.48"
$ws.Range("D3").Value = "1.652.03"
$ws.Range("D4").Value = "1.008"
$ws.Range("D5").Value = "312.11"
$ws.Range("D6").Value = "1.006"
$ws.Range("D7").Value = "0.3919"
$ws.Range("D8").Value = "0.3900"
$ws.Range("D9").Value = "1.008"
$ws.Range("D10").Value = "50.54"
$ws.Range("D12").Value = "0.08540"
$ws.Range("D13").Value = "24.95"
$ws.Range("D14").Value = "7.225"
$ws.Range("D15").Value = "0.00001302"
$ws.Range("D16").Value = "7.592"
$ws.Range("D17").Value = "1.669.52"
$ws.Range("D18").Value = "93.12"
$ws.Range("D19").Value = "0.06976"
$ws.Range("D20").Value = "21.28"
$ws.Range("D21").Value = "6.986"
$ws.Range("D22").Value = "1.007"
$ws.Range("D23").Value = "13.79"
$ws.Range("D24").Value = "24.392.25"
$ws.Range("D25").Value = "2.340"
$ws.Range("D26").Value = "2.760"
$ws.Range("D27").Value = "22.74"
$ws.Range("D28").Value = "5.785"
$ws.Range("D29").Value = "158.50"
$ws.Range("D30").Value = "144.63"
$ws.Range("D31").Value = "8.266"
$ws.Range("D32").Value = "2.517"
$ws.Range("D33").Value = "1.851.28"
$ws.Range("D34").Value = "0.08174"
$ws.Range("D35").Value = "0.03019"
$ws.Range("D36").Value = "0.9964"
$ws.Range("D37").Value = "6.859"
$ws.Range("D38").Value = "0.2767"
$ws.Range("D39").Value = "0.09543"
$ws.Range("D40").Value = "1.500"
$ws.Range("D41").Value = "10.19"
$ws.Range("D42").Value = "0.7772"
$ws.Range("D43").Value = "13.28"
$ws.Range("D44").Value = "16.33"
$ws.Range("D45").Value = "2.541"
$ws.Range("D47").Value = "4.149"
$ws.Range("D48").Value = "1.005"
$ws.Range("D49").Value = "0.08565"
$ws.Range("D50").Value = "1.311"
$ws.Range("D51").Value = "136.85"

# Update Volume(1h) (E) column
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("E3").Value = "  -3.20%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  -2.25%  "
$ws.Range("E8").Value = "  -3.48%  "
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  -5.75%  "
$ws.Range("E11").Value = "  -6.26%  "
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("E13").Value = "  -5.11%  "
$ws.Range("E14").Value = "  -3.89%  "
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("E16").Value = "  -5.09%  "
$ws.Range("E17").Value = "  -3.84%  "
$ws.Range("E18").Value = "  -2.63%  "
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("E21").Value = "  -4.26%  "
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  -4.58%  "
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").Value = "  -4.67%  "
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  -6.77%  "
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("E32").Value = "  +10.54%  "
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("E34").Value = "  -5.53%  "
$ws.Range("E35").Value = "  -5.61%  "
$ws.Range("E36").Value = "  -3.31%  "
$ws.Range("E37").Value = "  -6.06%  "
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("E41").Value = "  -4.59%  "
$ws.Range("E42").Value = "  -7.59%  "
$ws.Range("E43").Value = "  -6.73%  "
$ws.Range("E44").Value = "  -6.72%  "
$ws.Range("E45").Value = "  -6.71%  "
$ws.Range("E46").Value = "  -5.67%  "
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("E51").Value = "  -2.58%  "
